$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Columns.Item(2).Insert()

# New column B width should match column A's width
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Header row
$ws.Range("B1").Value = "StatQuery"

# New query cell in row 2, with wrap text like A2
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Colorectal cancer, NOS'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true
